$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '307.92'
Set-TextValue 'E2' '1.13%'
Set-TextValue 'D3' '36.36'
Set-TextValue 'E3' '1.43%'
Set-TextValue 'D4' '5.059'
Set-TextValue 'E4' '-0.27%'
Set-TextValue 'E5' '0.59%'
Set-TextValue 'D6' '2.104'
Set-TextValue 'E6' '9.28%'
Set-TextValue 'B7' 'GateToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.154'
Set-TextValue 'E7' '-0.17%'
Set-TextValue 'B8' 'KuCoinToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D8' '7.852'
Set-TextValue 'E8' '0.29%'
Set-TextValue 'B9' 'MXToken'
Set-TextValue 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9300'
Set-TextValue 'E9' '-0.29%'
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1461'
Set-TextValue 'E10' '15.71%'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1928'
Set-TextValue 'E11' '0.80%'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.09108'
Set-TextValue 'E12' '-1.31%'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03450'
Set-TextValue 'E13' '-0.80%'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09911'
Set-TextValue 'E14' '0.48%'
Set-TextValue 'B15' 'BitForexToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001407'
Set-TextValue 'E15' '-0.76%'
Set-TextValue 'B16' 'TigerCash'
Set-TextValue 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D16' '0.006311'
Set-TextValue 'E16' '-5.16%'
Set-TextValue 'B17' 'LEO'
Set-TextValue 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.848'
Set-TextValue 'E17' '6.39%'
Set-TextValue 'D18' '3.375'
Set-TextValue 'E18' '4.50%'
Set-TextValue 'D19' '0.3462'
Set-TextValue 'E19' '1.15%'
Set-TextValue 'D20' '0.1281'
Set-TextValue 'E20' '-4.10%'
Set-TextValue 'D21' '4.801'
Set-TextValue 'E21' '-7.27%'
Set-TextValue 'D22' '0.2339'
Set-TextValue 'D23' '0.04388'
Set-TextValue 'E23' '-0.22%'
Set-TextValue 'D24' '0.001232'
Set-TextValue 'E24' '-0.40%'
Set-TextValue 'E25' '4.31%'
Set-TextValue 'D27' '0.0001301'
Set-TextValue 'E27' '-0.13%'
Set-TextValue 'D39' '0.02021'
Set-TextValue 'E39' '2.18%'
Set-TextValue 'D40' '0.05171'
Set-TextValue 'E40' '0.28%'
Set-TextValue 'D41' '0.007486'
Set-TextValue 'E41' '-0.78%'
Set-TextValue 'D42' '0.01014'
Set-TextValue 'E42' '1.05%'
Set-TextValue 'D43' '0.1368'
Set-TextValue 'E43' '0.16%'
Set-TextValue 'D44' '0.002172'
Set-TextValue 'E44' '3.19%'
Set-TextValue 'D45' '0.009954'
Set-TextValue 'E45' '-6.75%'
Set-TextValue 'D46' '0.00006286'
Set-TextValue 'E46' '-0.83%'
Set-TextValue 'D47' '0.00000000749'
Set-TextValue 'E47' '-0.27%'
Set-TextValue 'D48' '64.85'
Set-TextValue 'E48' '-0.16%'
Set-TextValue 'D49' '0.001249'
Set-TextValue 'E49' '-21.99%'
Set-TextValue 'D50' '0.00002098'
Set-TextValue 'E50' '-0.27%'
Set-TextValue 'D51' '0.0001998'
Set-TextValue 'E51' '-0.27%'
